# "time update and bug fixing"
#
# - sequences (sheet1): selection changes to A1:E1 (was B1 / B:B).
# - Sheet3: loses tabSelected (it's no longer the active tab).
# - A brand-new sheet named "Sheet1" (4th tab, sheetId 4) is appended after
#   Sheet3, becomes the active sheet/tab (activeTab=3), and is populated with
#   an 11-row Point/Label/A/B/C table (same header layout as "subsequences"),
#   with selection G15.

$wb = $excel.ActiveWorkbook

# --- sequences: update the saved selection -------------------------------
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("A1:E1").Select()

# --- add the new trailing sheet ------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# Header row (reuses the same shared strings as the other sheets: Point,
# Label, A, B, C).
$ws.Range("A1").Value = "Point"
$ws.Range("B1").Value = "Label"
$ws.Range("C1").Value = "A"
$ws.Range("D1").Value = "B"
$ws.Range("E1").Value = "C"

# Data rows: Point, Label, A, B, C (one-hot of the label).
$data = @(
    @(1, "A", 1, 0, 0),
    @(2, "A", 1, 0, 0),
    @(3, "B", 0, 1, 0),
    @(4, "C", 0, 0, 1),
    @(5, "A", 1, 0, 0),
    @(6, "B", 0, 1, 0),
    @(7, "C", 0, 0, 1),
    @(8, "B", 0, 1, 0),
    @(9, "A", 1, 0, 0),
    @(10, "B", 0, 1, 0),
    @(11, "C", 0, 0, 1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Center-align the table, matching the style used on the other sheets.
$ws.Range("A1:E12").HorizontalAlignment = -4108

# Saved selection on the new sheet.
[void]$ws.Range("G15").Select()
